# Update the representative paragraph:
# "COFINA Gabon est représentée par Monsieur El Hadji Mamadou FAYE, son Directeur Général,"
# becomes
# "COFINA Gabon est représentée par Madame Jenny MVOU, son Directeur Général Adjointe,"

$d = $word.ActiveDocument

# Step A: "Monsieur" -> "Madame" (non-bold run, text stays non-bold)
$d.Content.Find.Execute("est représentée par Monsieur ", $true, $false, $false, $false, $false, $true, 1, $false, "est représentée par Madame ", 2)

# Step B: swap the bold name; including the trailing ", " pulls the comma/space into
# the bold run's formatting, matching the target "Jenny MVOU, " bold run.
$d.Content.Find.Execute("El Hadji Mamadou FAYE, ", $true, $false, $false, $false, $false, $true, 1, $false, "Jenny MVOU, ", 2)

# Step C: fix up the job title (all non-bold text, spans former runs "son Directeur " + "Généra" + "l,")
$d.Content.Find.Execute("son Directeur Général,", $true, $false, $false, $false, $false, $true, 1, $false, "son Directeur Général Adjointe,", 2)
